# Apply the Monster Hunter review edits.
#
# Plain Find/Replace in this runtime rebuilds the touched paragraph's runs
# and silently drops any *empty* <w:r/> run that sits next to the edited
# run once both share identical (empty) formatting - but the target XML
# keeps those empty runs untouched. To stay faithful to the original
# markup we instead locate the whole paragraph whose text equals the old
# string and surgically replace only the trailing run content (leaving
# <w:pPr> and any preceding empty run exactly as they were) via
# Range.InsertXML.

$d = $word.ActiveDocument

function Set-ParagraphTail($oldText, $innerRunXml) {
    foreach ($p in $d.Paragraphs) {
        $full = $p.Range.Text
        $trimmed = $full.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $oldText) {
            $r = $p.Range
            $tailLen = $oldText.Length
            $target = $d.Range($r.End - 1 - $tailLen, $r.End - 1)
            $xml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" `
                 + "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" `
                 + "<pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" `
                 + "<w:body><w:p>" + $innerRunXml + "</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
            $target.InsertXML($xml)
            return $true
        }
    }
    return $false
}

# Title (H1 heading - single plain run, no leading empty run)
Set-ParagraphTail `
    "Play Monster Hunter for Free - Impressive Graphics & Exciting Battle Mode" `
    '<w:r><w:t>Play Monster Hunter Free - Exciting and Immersive Slot Game</w:t></w:r>' | Out-Null

# "What we like" bullet list (each paragraph starts with an empty <w:r/>)
Set-ParagraphTail "Impressive graphics and design" `
    '<w:r><w:t>Impressive Graphics and Design</w:t></w:r>' | Out-Null
Set-ParagraphTail "Exciting and immersive battle mode" `
    '<w:r><w:t>Symbols and Functions</w:t></w:r>' | Out-Null
Set-ParagraphTail "Variety of character abilities" `
    '<w:r><w:t>Battle Mode</w:t></w:r>' | Out-Null
Set-ParagraphTail "Strategic gameplay potential" `
    '<w:r><w:t>Character Abilities</w:t></w:r>' | Out-Null

# "What we don't like" bullet list
Set-ParagraphTail "Symbols not as numerous as other slots" `
    '<w:r><w:t>Limited symbol variety in base game</w:t></w:r>' | Out-Null
Set-ParagraphTail "May require time to master gameplay strategy" `
    '<w:r><w:t>Requires strategic gameplay</w:t></w:r>' | Out-Null

# Closing bold title repeat
Set-ParagraphTail `
    "Play Monster Hunter for Free - Impressive Graphics & Exciting Battle Mode" `
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Play Monster Hunter Free - Exciting and Immersive Slot Game</w:t></w:r>' | Out-Null

# Closing italic meta description
Set-ParagraphTail `
    "Experience the unique and immersive gameplay of Monster Hunter with exceptional graphics, character abilities, and strategic potential. Play for free now!" `
    '<w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Monster Hunter, a free online slot game with impressive graphics and strategic gameplay.</w:t></w:r>' | Out-Null
